$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# --- Insert rows for the newly-waited indicators that were added to the union ---
# 4 new rows for INDICATOR_101, INDICATOR_102, INDICATOR_103, INDICATOR_105
$ws.Rows("53:56").Insert()
# 2 new rows for INDICATOR_112, INDICATOR_113
$ws.Rows("96:97").Insert()

# --- Make sure the newly-inserted blank rows carry the same A/B/E/F boilerplate
#     values as every other row in this block (Insert() only copies formatting) ---
for ($r = 53; $r -le 102; $r++) {
    $ws.Cells.Item($r, 1).Value2 = "CREATE/MODIFY"
    $ws.Cells.Item($r, 2).Value2 = "LIB_EWS_RETAIL"
    $ws.Cells.Item($r, 5).Value2 = "String"
    $ws.Cells.Item($r, 6).Value2 = "String"
}

# --- Rewrite column C (Key) for rows 53-102 so the union list is in its new order ---
$ws.Cells.Item(53, 3).Value2 = "INDICATOR_101"
$ws.Cells.Item(54, 3).Value2 = "INDICATOR_102"
$ws.Cells.Item(55, 3).Value2 = "INDICATOR_103"
$ws.Cells.Item(56, 3).Value2 = "INDICATOR_105"
$ws.Cells.Item(57, 3).Value2 = "INDICATOR_114"
$ws.Cells.Item(58, 3).Value2 = "INDICATOR_200"
$ws.Cells.Item(59, 3).Value2 = "INDICATOR_201"
$ws.Cells.Item(60, 3).Value2 = "INDICATOR_202"
$ws.Cells.Item(61, 3).Value2 = "INDICATOR_203"
$ws.Cells.Item(62, 3).Value2 = "INDICATOR_204"
$ws.Cells.Item(63, 3).Value2 = "INDICATOR_205"
$ws.Cells.Item(64, 3).Value2 = "INDICATOR_206"
$ws.Cells.Item(65, 3).Value2 = "INDICATOR_207"
$ws.Cells.Item(66, 3).Value2 = "INDICATOR_208"
$ws.Cells.Item(67, 3).Value2 = "INDICATOR_209"
$ws.Cells.Item(68, 3).Value2 = "INDICATOR_210"
$ws.Cells.Item(69, 3).Value2 = "INDICATOR_211"
$ws.Cells.Item(70, 3).Value2 = "INDICATOR_212"
$ws.Cells.Item(71, 3).Value2 = "INDICATOR_213"
$ws.Cells.Item(72, 3).Value2 = "INDICATOR_214"
$ws.Cells.Item(73, 3).Value2 = "INDICATOR_215"
$ws.Cells.Item(74, 3).Value2 = "INDICATOR_216"
$ws.Cells.Item(75, 3).Value2 = "INDICATOR_217"
$ws.Cells.Item(76, 3).Value2 = "INDICATOR_218"
$ws.Cells.Item(77, 3).Value2 = "INDICATOR_219"
$ws.Cells.Item(78, 3).Value2 = "INDICATOR_220"
$ws.Cells.Item(79, 3).Value2 = "INDICATOR_221"
$ws.Cells.Item(80, 3).Value2 = "INDICATOR_222"
$ws.Cells.Item(81, 3).Value2 = "INDICATOR_223"
$ws.Cells.Item(82, 3).Value2 = "INDICATOR_224"
$ws.Cells.Item(83, 3).Value2 = "INDICATOR_225"
$ws.Cells.Item(84, 3).Value2 = "INDICATOR_226"
$ws.Cells.Item(85, 3).Value2 = "INDICATOR_227"
$ws.Cells.Item(86, 3).Value2 = "INDICATOR_228"
$ws.Cells.Item(87, 3).Value2 = "INDICATOR_229"
$ws.Cells.Item(88, 3).Value2 = "INDICATOR_230"
$ws.Cells.Item(89, 3).Value2 = "INDICATOR_231"
$ws.Cells.Item(90, 3).Value2 = "INDICATOR_106"
$ws.Cells.Item(91, 3).Value2 = "INDICATOR_107"
$ws.Cells.Item(92, 3).Value2 = "INDICATOR_108"
$ws.Cells.Item(93, 3).Value2 = "INDICATOR_109"
$ws.Cells.Item(94, 3).Value2 = "INDICATOR_110"
$ws.Cells.Item(95, 3).Value2 = "INDICATOR_111"
$ws.Cells.Item(96, 3).Value2 = "INDICATOR_112"
$ws.Cells.Item(97, 3).Value2 = "INDICATOR_113"
$ws.Cells.Item(98, 3).Value2 = "INDICATOR_115"
$ws.Cells.Item(99, 3).Value2 = "INDICATOR_116"
$ws.Cells.Item(100, 3).Value2 = "INDICATOR_117"
$ws.Cells.Item(101, 3).Value2 = "INDICATOR_150"
$ws.Cells.Item(102, 3).Value2 = "INDICATOR_151"

# --- Fix up the two rows whose cell style differs from their neighbours ---
$ws.Range("C2").Copy()
$ws.Range("C56").PasteSpecial(-4122)
$ws.Range("C97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view: scroll target + single-cell selection on C94 ---
$ws.Activate()
$ws.Range("C94").Select()
